$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CDRDfRCP")
$ws.Range("B3").Value = 0.05
